$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Files-count query (same text used in C2, C3, C4)
$filesCountQuery = @'
MATCH (f:file)
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina NovaSeq']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,f, s, collect(distinct samp.sample_id) as samp
RETURN
count(distinct s) AS Studies,
count(distinct p) AS Participants,
count(distinct samp) AS Samples,
count(distinct f) AS Files
'@

# Participants query (row 2, column B) - note the trailing space after ['Illumina NovaSeq']
$participantsQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina NovaSeq'] 
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p, s, collect(distinct samp.sample_id) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY `Participant ID`LIMIT 100
'@

# Samples query (row 3, column B)
$samplesQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina NovaSeq']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as `Sample ID`,
 coalesce(p.participant_id,'') as `Participant ID`,
 coalesce(s.study_name, '') as `Study Name`,
 coalesce(s.phs_accession,'') as `Accession`,
coalesce(samp.sample_tumor_status,'') as `Tumor`,
coalesce(samp.sample_type,'') as `Analyte Type`
ORDER By samp.sample_id LIMIT 100
'@

# Files query (row 4, column B)
$filesQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina NovaSeq']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,f,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN 
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name, '') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id,'') as `Participant ID`,
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
ORDER By f.file_name LIMIT 100
'@

$neo4jFile = "TC10_CDS_Filter_InstrumentModel-Illumina NovaSeq_Neo4jData.xlsx"
$webFile = "TC10_CDS_Filter_InstrumentModel-Illumina NovaSeq_WebData.xlsx"

# Write order matters for the shared-string table layout: filenames first,
# then the B-column queries (Participants, Samples, Files), then the
# C-column "files count" query last (it is shared across C2:C4).
$ws.Range("D2").Value = $neo4jFile
$ws.Range("E2").Value = $webFile
$ws.Range("D3").Value = $neo4jFile
$ws.Range("E3").Value = $webFile
$ws.Range("D4").Value = $neo4jFile
$ws.Range("E4").Value = $webFile

$ws.Range("B2").Value = $participantsQuery
$ws.Range("B3").Value = $samplesQuery
$ws.Range("B4").Value = $filesQuery

$ws.Range("C2").Value = $filesCountQuery
$ws.Range("C3").Value = $filesCountQuery
$ws.Range("C4").Value = $filesCountQuery

# Column D/E widened to fit the new (longer) file-name text
$ws.Columns.Item(4).ColumnWidth = 91.25
$ws.Columns.Item(5).ColumnWidth = 89.6

# Selection moves to D2 when the file was last saved
$ws.Range("D2").Select()
